$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "How many curve shades can I create?"
$ws.Range("B40").Value = "llama3.2:latest"
$ws.Range("C40").Value = "According to the document, you can create 250 curve shades."

$ws.Range("A41").Value = "What's the maximum number of lithology types in an log?"
$ws.Range("B41").Value = "llama3.2:latest"
$ws.Range("C41").Value = "The maximum number of lithology types that can be recorded in a log is 450."

$ws.Range("A42").Value = "How many tracks can you define in one ODF?"
$ws.Range("B42").Value = "llama3.2:latest"
$ws.Range("C42").Value = "According to the Track Settings dialog box, the number of tracks that can be defined is 200."
